$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with data about NOAA OI SST V2 High Resolution Dataset
$ws.Range("A5").Value = "NOAA OI SST V2 High Resolution Dataset"
$ws.Range("B5").Value = "0.25 degree"
$ws.Range("C5").Value = "daily"
$ws.Range("D5").Value = "1981 – present"
$ws.Range("E5").Value = "NOAA"
$ws.Range("F5").Value = "no location"
$ws.Range("G5").Value = "netCDF"
$ws.Range("H5").Value = "https://psl.noaa.gov/data/gridded/data.noaa.oisst.v2.highres.html#detail"

# Update the hyperlink display text on H4 to match the cell text.
# (Mutating TextToDisplay on an existing Hyperlink object in-place isn't
# supported by this host - it silently creates a duplicate hyperlink - so
# remove the old one and re-add it with the same address/target and the
# new display text.)
$ws.Range("H4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H4"), "https://polar.ncep.noaa.gov/sst/ophi/", "", "", "https://polar.ncep.noaa.gov/sst/ophi/ no longer available?")

# Move the selection to D11
$ws.Range("D11").Select()
